$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text in the source workbook (inline
# strings), even for values that look numeric (e.g. "310.63"). Assigning
# such a string straight to .Value would make Excel auto-convert it to a
# real number, so for those cells we temporarily force a Text number format,
# write the literal string, then restore the original "Normal" style so no
# stray formatting is left behind.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D14", "D16", "D17", "D21", "D22", "D23", "D27", "D30", "D31", "D33", "D34", "D35", "D36", "D40", "D41", "D43", "D45", "D46", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.496.67"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.468.07"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "310.63"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "91.48"
$ws.Range("E6").Value = "  -3.50%  "
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -4.32%  "
$ws.Range("D10").Value = "31.72"
$ws.Range("E10").Value = "  -6.37%  "
$ws.Range("D11").Value = "0.0769"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "2.846.25"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "6.71"
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("D15").Value = "2.458.32"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "15.07"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "0.753"
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("D18").Value = "41.401.94"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "70.19"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "10.91"
$ws.Range("E22").Value = "  -5.55%  "
$ws.Range("D23").Value = "233.09"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("D27").Value = "24.03"
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").Value = "35.82"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "153.15"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -5.58%  "
$ws.Range("D33").Value = "2.55"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").Value = "0.0749"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "17.68"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").Value = "2.48"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("D40").Value = "0.0995"
$ws.Range("E40").Value = "  -5.76%  "
$ws.Range("D41").Value = "4.07"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "19.23"
$ws.Range("E43").Value = "  -9.25%  "
$ws.Range("D44").Value = "1.939.48"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "2.711.43"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "94.67"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("E50").Value = "  -4.77%  "
$ws.Range("D51").Value = "65.55"
$ws.Range("E51").Value = "  -5.70%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
